# Fix Excel export by removing timezone from datetime columns.
#
# 1) Column B ("updated_at") previously carried two distinct timestamps
#    that differed only by the time-of-day component (i.e. a timezone
#    offset baked into the export). Collapse every row onto a single,
#    timezone-free timestamp.
# 2) Rows 7-9 (parts NP004/NP005/NP006) pick up the rest of their record
#    now that the PO-sourced fields have been merged in (cost, currency,
#    remarks, PO number, quantity, order/delivery dates, vendor, and a
#    cleaned-up "sources" provenance blob).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = 46028.67361742735

# --- 1) Normalize updated_at (column B) for every data row ---------------
foreach ($r in 2..11) {
    $ws.Range("B$r").Value = $newTimestamp
}

# --- helper: write a value as TEXT (shared string) without leaving the  --
# --- cell's style index changed (matches the plain, unstyled text cells -
# --- used throughout the rest of the sheet) -------------------------------
function Set-TextValue {
    param($range, $text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$poSources = '[{"source_system": "pos", "source_file": "po_mock_newparts.pdf"}, {"source_system": "user", "source_file": "diagram2.pdf"}]'

# --- 2) Fill in the newly-merged PO fields for rows 7-9 -------------------
Set-TextValue $ws.Range("J7") "1683.78"
Set-TextValue $ws.Range("AA7") "INR"
Set-TextValue $ws.Range("AD7") "Standard Order"
Set-TextValue $ws.Range("AJ7") "PO1003"
Set-TextValue $ws.Range("AM7") "471"
Set-TextValue $ws.Range("AV7") "2024-02-10"
Set-TextValue $ws.Range("AW7") "2024-04-11"
Set-TextValue $ws.Range("BD7") "MegaTools Pvt Ltd"
Set-TextValue $ws.Range("BB7") $poSources

Set-TextValue $ws.Range("J8") "1428.22"
Set-TextValue $ws.Range("AA8") "USD"
Set-TextValue $ws.Range("AD8") "Check Material Cert"
Set-TextValue $ws.Range("AJ8") "PO1004"
Set-TextValue $ws.Range("AM8") "44"
Set-TextValue $ws.Range("AV8") "2024-11-13"
Set-TextValue $ws.Range("AW8") "2024-09-13"
Set-TextValue $ws.Range("BD8") "ElectroMart"
Set-TextValue $ws.Range("BB8") $poSources

Set-TextValue $ws.Range("J9") "16.41"
Set-TextValue $ws.Range("AA9") "USD"
Set-TextValue $ws.Range("AD9") "Urgent"
Set-TextValue $ws.Range("AJ9") "PO1005"
Set-TextValue $ws.Range("AM9") "445"
Set-TextValue $ws.Range("AV9") "2024-11-08"
Set-TextValue $ws.Range("AW9") "2024-02-25"
Set-TextValue $ws.Range("BD9") "UniMach"
Set-TextValue $ws.Range("BB9") $poSources
